# Update the "sections" worksheet to better represent the data model:
#  - rename the "localisation" / "industrie" headers to "country" / "region"
#  - change the sample "Paris" location values to "France"
#  - add eight new trailing columns (product_type_1-3, currency,
#    line_of_business, industry, sic_code, include) with their headers

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sections")

# Rename existing headers.
$ws.Range("E1").Value = "country"
$ws.Range("F1").Value = "region"

# Update the existing "Paris" values to "France".
$ws.Range("E3").Value = "France"
$ws.Range("E4").Value = "France"

# Insert eight new columns after F (they become G:N) so the existing
# formatting on column F is not disturbed, then populate their headers.
$ws.Range("G1:N1").EntireColumn.Insert()

$ws.Range("G1").Value = "product_type_1"
$ws.Range("H1").Value = "product_type_2"
$ws.Range("I1").Value = "product_type_3"
$ws.Range("J1").Value = "currency"
$ws.Range("K1").Value = "line_of_business"
$ws.Range("L1").Value = "industry"
$ws.Range("M1").Value = "sic_code"
$ws.Range("N1").Value = "include"
